# students.xlsx update
#  - flips a batch of IsSocialCase (J) / IsMedicalCase (K) cells that were
#    "TRUE" over to FALSE (some become native booleans, some stay text
#    "FALSE", matching whatever the source row already looked like)
#  - swaps the data that was entered on the wrong rows (81 <-> 82)
#  - re-autofits col B, and leaves the selection on F62 having scrolled
#    the window up a bit

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- J (IsSocialCase) / K (IsMedicalCase) corrections ------------------

$J_BOOL = @(4,7,8,9,10)
$J_TEXT = @(13,14,16,23,24,28,33,34,35,36,37,41,42,43,45,46,48,54,55,56,60,61,62,63,65,66,69,73,74,75,78,79)
$K_BOOL = @(2,7,8,10,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,71,72,73,74,75,76,77,78,79,80)
$K_TEXT = @(16,17,18,19,20,22,23,24,25,26,29,31,32,34,36,39,40,41,42)

foreach ($r in $J_BOOL) {
    $ws.Cells.Item($r, 10).Value = $false
}
foreach ($r in $K_BOOL) {
    $ws.Cells.Item($r, 11).Value = $false
}

# these need to stay literal text "FALSE" (not a boolean) - force text entry
# with a leading apostrophe, then strip the resulting quote-prefix format so
# the cell is left with no explicit style, same as its neighbours.
foreach ($r in $J_TEXT) {
    $cell = $ws.Cells.Item($r, 10)
    $cell.Value = "'FALSE"
    $cell.ClearFormats()
}
foreach ($r in $K_TEXT) {
    $cell = $ws.Cells.Item($r, 11)
    $cell.Value = "'FALSE"
    $cell.ClearFormats()
}

# ---- rows 81 / 82 had been entered swapped - fix them ------------------

$ws.Range("C81").Value = 456382
$ws.Range("D81").Value = 2123456789123
$ws.Range("E81").Value = "Beatrice"
$ws.Range("F81").Value = "Vaduva"
$ws.Range("G81").Value = "D"
$ws.Range("Q81").Value = 743456789
$ws.Range("R81").Value = 80

$ws.Range("C82").Value = "'TC 419786"
$ws.Range("C82").ClearFormats()
$ws.Range("D82").Value = 2971201360023
$ws.Range("E82").Value = "Eliza Ioana"
$ws.Range("F82").Value = "Țuțuianu"
$ws.Range("G82").Value = "L"
$ws.Range("Q82").Value = 748318768
$ws.Range("R82").ClearContents()

# ---- cosmetics: column B autofit + view position ------------------------

$ws.Columns("B").AutoFit()

$ws.Range("F62").Select()
$excel.ActiveWindow.ScrollRow = 55
